$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the example data in the table (row 2-5) with more realistic values
$ws.Range("B2").Value = "SSH1050"
$ws.Range("D2").Value = "GV1234"
$ws.Range("E2").Value = "20161234, 201612345, 20161236, 20161237"

$ws.Range("B3").Value = "EM1170"
$ws.Range("D3").Value = "GV1235"
$ws.Range("E3").Value = "20161234, 201612345, 20161236, 20161237"

$ws.Range("B4").Value = "IT1110"
$ws.Range("D4").Value = "GV1236"
$ws.Range("E4").Value = "20161234, 201612345, 20161236, 20161237"

$ws.Range("B5").Value = "MI1110"
$ws.Range("D5").Value = "GV1237"
$ws.Range("E5").Value = "20161234, 201612345, 20161236, 20161237"

# Move the active selection to D6, matching the saved selection state
$ws.Range("D6").Select()
